$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Charles Bediako"
$ws.Range("B31").Value = "Alabama"
$ws.Range("C31").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4565203.png&w=350&h=254"

$ws.Range("A32").Value = "Jaden Bradley"
$ws.Range("B32").Value = "Alabama"
$ws.Range("C32").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432737.png&w=350&h=254"

$ws.Range("A33").Value = "Nimari Burnett"
$ws.Range("B33").Value = "Alabama"
$ws.Range("C33").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4708027.png&w=350&h=254"

$ws.Range("A34").Value = "Noah Clowney"
$ws.Range("B34").Value = "Alabama"
$ws.Range("C34").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4712896.png&w=350&h=254"

$ws.Range("A35").Value = "Adam Cottrell"
$ws.Range("B35").Value = "Alabama"
$ws.Range("C35").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4598102.png&w=350&h=254"

$ws.Range("A36").Value = "Rylan Griffen"
$ws.Range("B36").Value = "Alabama"
$ws.Range("C36").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4683682.png&w=350&h=254"

$ws.Range("A37").Value = "Noah Gurley"
$ws.Range("B37").Value = "Alabama"
$ws.Range("C37").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4280015.png&w=350&h=254"

$ws.Range("A38").Value = "Delaney Heard"
$ws.Range("B38").Value = "Alabama"
$ws.Range("C38").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4702175.png&w=350&h=254"

$ws.Range("A39").Value = "Darius Miles"
$ws.Range("B39").Value = "Alabama"
$ws.Range("C39").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4433561.png&w=350&h=254"

$ws.Range("A40").Value = "Brandon Miller"
$ws.Range("B40").Value = "Alabama"
$ws.Range("C40").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4433287.png&w=350&h=254"

$ws.Range("A41").Value = "Nick Pringle"
$ws.Range("B41").Value = "Alabama"
$ws.Range("C41").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4703887.png&w=350&h=254"

$ws.Range("A42").Value = "Jaden Quinerly"
$ws.Range("B42").Value = "Alabama"
$ws.Range("C42").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4702179.png&w=350&h=254"

$ws.Range("A43").Value = "Jahvon Quinerly"
$ws.Range("B43").Value = "Alabama"
$ws.Range("C43").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4397132.png&w=350&h=254"

$ws.Range("A44").Value = "Max Scharnowski"
$ws.Range("B44").Value = "Alabama"
$ws.Range("C44").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105513.png&w=350&h=254"

$ws.Range("A45").Value = "Mark Sears"
$ws.Range("B45").Value = "Alabama"
$ws.Range("C45").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4703530.png&w=350&h=254"

$ws.Range("A46").Value = "Kai Spears"
$ws.Range("B46").Value = "Alabama"
$ws.Range("C46").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105514.png&w=350&h=254"

$ws.Range("A47").Value = "Dominick Welch"
$ws.Range("B47").Value = "Alabama"
$ws.Range("C47").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4397221.png&w=350&h=254"

$ws.Range("A48").Value = "Jaden Akins"
$ws.Range("B48").Value = "Michigan State"
$ws.Range("C48").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4683730.png&w=350&h=254"

$ws.Range("A49").Value = "Pierre Brooks"
$ws.Range("B49").Value = "Michigan State"
$ws.Range("C49").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4683731.png&w=350&h=254"

$ws.Range("A50").Value = "Carson Cooper"
$ws.Range("B50").Value = "Michigan State"
$ws.Range("C50").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105817.png&w=350&h=254"

$ws.Range("A51").Value = "Malik Hall"
$ws.Range("B51").Value = "Michigan State"
$ws.Range("C51").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4592693.png&w=350&h=254"

$ws.Range("A52").Value = "Joey Hauser"
$ws.Range("B52").Value = "Michigan State"
$ws.Range("C52").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4295180.png&w=350&h=254"

$ws.Range("A53").Value = "AJ Hoggard"
$ws.Range("B53").Value = "Michigan State"
$ws.Range("C53").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432206.png&w=350&h=254"

$ws.Range("A54").Value = "Tre Holloman"
$ws.Range("B54").Value = "Michigan State"
$ws.Range("C54").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105816.png&w=350&h=254"

$ws.Range("A55").Value = "Steven Izzo"
$ws.Range("B55").Value = "Michigan State"
$ws.Range("C55").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4592692.png&w=350&h=254"

$ws.Range("A56").Value = "Jaxon Kohler"
$ws.Range("B56").Value = "Michigan State"
$ws.Range("C56").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105815.png&w=350&h=254"

$ws.Range("A57").Value = "Nick Sanders"
$ws.Range("B57").Value = "Michigan State"
$ws.Range("C57").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105818.png&w=350&h=254"

$ws.Range("A58").Value = "Mady Sissoko"
$ws.Range("B58").Value = "Michigan State"
$ws.Range("C58").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4433154.png&w=350&h=254"

$ws.Range("A59").Value = "Davis Smith"
$ws.Range("B59").Value = "Michigan State"
$ws.Range("C59").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4702244.png&w=350&h=254"

$ws.Range("A60").Value = "Tyson Walker"
$ws.Range("B60").Value = "Michigan State"
$ws.Range("C60").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432129.png&w=350&h=254"

$ws.Range("A61").Value = "Jason Whitens"
$ws.Range("B61").Value = "Michigan State"
$ws.Range("C61").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4279457.png&w=350&h=254"

$ws.Range("B58").Select()
